# Add a new bulleted list item "Extract some repeating code in Renderer.cs"
# right after the existing last item "Some constants extracted in Renderer.cs",
# matching the same list formatting (pStyle a3 / numPr ilvl0 numId1) and moving
# the trailing "_GoBack" bookmark onto the new (now last) paragraph, exactly as
# Word itself would when the user types a new bullet at the end of the document.

$d = $word.ActiveDocument

# "_GoBack" marks the point of the most recent edit; it currently sits right
# after "Some constants extracted in Renderer.cs". Anchoring on the bookmark's
# own Range (rather than a freshly derived Range/position) lets the engine
# correctly grow/shift the bookmark together with the text we insert.
$bm = $d.Bookmarks("_GoBack")
$splitPos = $bm.Start

$r = $bm.Range
$r.InsertAfter("Extract some repeating code in Renderer.cs")

# Split the (still single) paragraph into two right before the text we just
# typed, turning it into its own list item while leaving the bookmark - which
# tracked the inserted text - at the end of the new paragraph.
$breakRange = $d.Range($splitPos, $splitPos)
$breakRange.InsertParagraphAfter()
